$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.151.61"
$ws.Range("E2").Value = "  +2.04%  "

# Row 3
$ws.Range("D3").Value = "3.459.03"
$ws.Range("E3").Value = "  +1.39%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.30"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.65"
$ws.Range("E6").Value = "  +2.44%  "

# Row 7
$ws.Range("D7").Value = "3.459.33"
$ws.Range("E7").Value = "  +1.37%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  +0.77%  "

# Row 10
$ws.Range("E10").Value = "  +2.85%  "

# Row 11
$ws.Range("E11").Value = "  +1.19%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.405"
$ws.Range("E12").Value = "  +5.10%  "

# Row 13
$ws.Range("D13").Value = "4.052.70"
$ws.Range("E13").Value = "  +1.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.29"
$ws.Range("E14").Value = "  +3.29%  "

# Row 15
$ws.Range("E15").Value = "  +2.49%  "

# Row 16
$ws.Range("D16").Value = "3.469.47"
$ws.Range("E16").Value = "  +1.53%  "

# Row 17
$ws.Range("E17").Value = "  +1.21%  "

# Row 18
$ws.Range("D18").Value = "63.120.54"
$ws.Range("E18").Value = "  +1.96%  "

# Row 19
$ws.Range("E19").Value = "  +3.97%  "

# Row 20
$ws.Range("E20").Value = "  +3.38%  "

# Row 21
$ws.Range("E21").Value = "  +1.50%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.21"
$ws.Range("E22").Value = "  -0.55%  "

# Row 23
$ws.Range("E23").Value = "  +1.79%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.56"
$ws.Range("E24").Value = "  -0.48%  "

# Row 25
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("D26").Value = "3.607.11"
$ws.Range("E26").Value = "  +1.47%  "

# Row 27
$ws.Range("E27").Value = "  +1.17%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.181"
$ws.Range("E28").Value = "  -1.77%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.64"
$ws.Range("E29").Value = "  +2.62%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.993"
$ws.Range("E30").Value = "  -0.61%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.16"
$ws.Range("E31").Value = "  +1.98%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.13"
$ws.Range("E32").Value = "  -0.15%  "

# Row 33
$ws.Range("E33").Value = "  +0.01%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.41"
$ws.Range("E34").Value = "  -0.56%  "

# Row 35
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.34"
$ws.Range("E35").Value = "  -3.84%  "

# Row 36
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.15"
$ws.Range("E36").Value = "  +2.58%  "

# Row 37
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.32"
$ws.Range("E37").Value = "  +1.60%  "

# Row 38
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "31.97"
$ws.Range("E38").Value = "  +12.09%  "

# Row 39
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.60"
$ws.Range("E39").Value = "  +4.31%  "

# Row 40
$ws.Range("E40").Value = "  +0.51%  "

# Row 41
$ws.Range("D41").Value = "3.495.66"
$ws.Range("E41").Value = "  +1.48%  "

# Row 42
$ws.Range("E42").Value = "  +2.07%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.792"
$ws.Range("E43").Value = "  +0.94%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.40"
$ws.Range("E44").Value = "  -0.75%  "

# Row 45
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.73"
$ws.Range("E45").Value = "  +3.54%  "

# Row 46
$ws.Range("E46").Value = "  +3.82%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.37"
$ws.Range("E47").Value = "  -1.26%  "

# Row 48
$ws.Range("D48").Value = "2.596.36"
$ws.Range("E48").Value = "  +3.79%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  +11.78%  "

# Row 50
$ws.Range("E50").Value = "  +2.86%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.96"
$ws.Range("E51").Value = "  +0.69%  "
